$p = $ppt.ActivePresentation

# Remove the "Homework" slide (slide 21). The following slide
# ("Free Trainings @ Telerik Academy") shifts up to become the new
# last slide (21) and its slide-number field updates automatically.
$p.Slides.Item(21).Delete()
